# Refresh the cryptocurrency Price (col D) and Volume(1h) (col E) figures to the
# latest scraped snapshot. Values that render as plain decimals (e.g. "534.41")
# must stay text, matching the sheet's existing inline-string cells, so we force
# the text format before writing and clear the format delta afterwards so no
# style index is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.176.68"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "2.518.12"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "534.41"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.07"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.562"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.89%  "
$ws.Range("D9").Value = "2.526.04"
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0991"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.161"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.40"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.351"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").Value = "2.963.32"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("E15").Value = "  -3.33%  "
$ws.Range("D16").Value = "59.136.48"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D18").Value = "2.502.68"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("E19").Value = "  -3.43%  "
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "319.48"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.78"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.87"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.418"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -5.27%  "
$ws.Range("E26").Value = "  +2.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.997"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.75"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.45%  "
$ws.Range("D30").Value = "0.0₃0764"
$ws.Range("E30").Value = "  -1.87%  "
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "162.88"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.26%  "
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("E34").Value = "  -9.20%  "
$ws.Range("E35").Value = "  -0.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.42"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.21"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -4.96%  "
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.87"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.35"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -8.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.63"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.41%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "287.42"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -6.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.801"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.25%  "
$ws.Range("E44").Value = "  +0.38%  "
$ws.Range("E45").Value = "  +0.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.597"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.61"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0926"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.54"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("E50").Value = "  -1.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0223"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.14%  "
